# Edit: (1) switch the presentation's table style on the slide-6 table to
# the built-in style {EA436D7E-2F65-4FB2-B298-D1948AE4A2B3}; (2) swap the
# deck's colour theme from "Integral" to "Office Theme" (the 12 theme
# colours that live in ppt/theme/theme1.xml, used by the slide master).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
$table  = $slide6.Shapes.Item(2).Table
$table.ApplyStyle("{EA436D7E-2F65-4FB2-B298-D1948AE4A2B3}", $false)

# --- 2. Theme colours: Integral -> Office Theme ---------------------------
# Colors() is 1-based and follows clrScheme document order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
# 9 accent5, 10 accent6, 11 hlink, 12 folHlink.
# RGB is packed little-endian (R + G*256 + B*65536), matching VBA's RGB().
$themeColors = $p.Slides.Item(1).ThemeColorScheme
$themeColors.Colors(1).RGB  = 0          # dk1      000000
$themeColors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$themeColors.Colors(3).RGB  = 6968388    # dk2      44546A
$themeColors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$themeColors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$themeColors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$themeColors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$themeColors.Colors(8).RGB  = 49407      # accent4  FFC000
$themeColors.Colors(9).RGB  = 12874308   # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456    # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797   # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477    # folHlink 954F72
